$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.751.18"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "2.376.64"
$ws.Range("E3").Value = "  -4.26%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "478.17"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.21"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").Value = "2.380.39"
$ws.Range("E9").Value = "  -4.75%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("E11").Value = "  -6.03%  "
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").Value = "2.789.75"
$ws.Range("E14").Value = "  -4.35%  "
$ws.Range("D15").Value = "55.867.62"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("E16").Value = "  -4.14%  "
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").Value = "2.380.73"
$ws.Range("E18").Value = "  -4.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.57"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "315.46"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.68"
$ws.Range("E21").Value = "  -4.90%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "56.84"
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.395"
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("E27").Value = "  -5.40%  "
$ws.Range("D28").Value = "2.482.47"
$ws.Range("E28").Value = "  -4.96%  "
$ws.Range("E29").Value = "  -6.02%  "
$ws.Range("D30").Value = "0.0₃0770"
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.80"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.01"
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  -3.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.57"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.833"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.41"
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("E43").Value = "  -3.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0948"
$ws.Range("E44").Value = "  +3.96%  "
$ws.Range("E45").Value = "  -5.33%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "255.71"
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.54"
$ws.Range("E49").Value = "  -6.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.92"
$ws.Range("E50").Value = "  -3.95%  "
$ws.Range("D51").Value = "1.777.05"
$ws.Range("E51").Value = "  -6.77%  "
